$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1653.9231
$ws.Range("I32").Value = 1833.3334
$ws.Range("J32").Value = 1600.1
$ws.Range("K32").Value = 1833.3334
$ws.Range("L32").Value = 1600.1
$ws.Range("M32").Value = -1507.3334
$ws.Range("N32").Value = -2252.1
$ws.Range("H121").Value = 3316
$ws.Range("I121").Value = 1661.6666
$ws.Range("J121").Value = 5797.5
$ws.Range("K121").Value = 4984.9998
$ws.Range("L121").Value = 17392.5
$ws.Range("M121").Value = -3237.9998
$ws.Range("N121").Value = -20886.5
$ws.Range("H137").Value = 2410.276
$ws.Range("I137").Value = 4452
$ws.Range("J137").Value = 1632.4762
$ws.Range("K137").Value = 13356
$ws.Range("L137").Value = 4897.4286
$ws.Range("M137").Value = -10806
$ws.Range("N137").Value = -9997.428599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4253.3184
$ws.Range("I61").Value = 2141.7273
$ws.Range("J61").Value = 6364.909
$ws.Range("K61").Value = 2141.7273
$ws.Range("L61").Value = 6364.909
$ws.Range("M61").Value = -1929.7273
$ws.Range("N61").Value = -6788.909
$ws.Range("H74").Value = 1480.4
$ws.Range("I74").Value = 1351.742
$ws.Range("J74").Value = 1765.2858
$ws.Range("K74").Value = 1351.742
$ws.Range("L74").Value = 1765.2858
$ws.Range("M74").Value = -477.742
$ws.Range("N74").Value = -3513.2858
$ws.Range("H77").Value = 1480.4
$ws.Range("I77").Value = 1351.742
$ws.Range("J77").Value = 1765.2858
$ws.Range("K77").Value = 6758.71
$ws.Range("L77").Value = 8826.429
$ws.Range("M77").Value = -2390.71
$ws.Range("N77").Value = -17562.429
$ws.Range("H132").Value = 6353.7856
$ws.Range("I132").Value = 3274.25
$ws.Range("J132").Value = 7585.6
$ws.Range("K132").Value = 9822.75
$ws.Range("L132").Value = 22756.8
$ws.Range("M132").Value = -7292.75
$ws.Range("N132").Value = -27816.8
$ws.Range("H136").Value = 4253.3184
$ws.Range("I136").Value = 2141.7273
$ws.Range("J136").Value = 6364.909
$ws.Range("K136").Value = 6425.1819
$ws.Range("L136").Value = 19094.727
$ws.Range("M136").Value = -3875.1819
$ws.Range("N136").Value = -24194.727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 970.3333
$ws.Range("I107").Value = 955.5
$ws.Range("K107").Value = 955.5
$ws.Range("M107").Value = 964.5
$ws.Range("H134").Value = 6806.6665
$ws.Range("I134").Value = 2788.9412
$ws.Range("J134").Value = 11075.5
$ws.Range("K134").Value = 8366.8236
$ws.Range("L134").Value = 33226.5
$ws.Range("M134").Value = -5831.8236
$ws.Range("N134").Value = -38296.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5849921
$ws.Range("I31").Value = 1572.2927
$ws.Range("J31").Value = 20836314
$ws.Range("K31").Value = 1572.2927
$ws.Range("L31").Value = 20836314
$ws.Range("M31").Value = -1277.2927
$ws.Range("N31").Value = -20836904
$ws.Range("H34").Value = 5849921
$ws.Range("I34").Value = 1572.2927
$ws.Range("J34").Value = 20836314
$ws.Range("K34").Value = 1572.2927
$ws.Range("L34").Value = 20836314
$ws.Range("M34").Value = -1370.2927
$ws.Range("N34").Value = -20836718
$ws.Range("H58").Value = 2967.2646
$ws.Range("I58").Value = 1832.2307
$ws.Range("J58").Value = 6656.125
$ws.Range("K58").Value = 1832.2307
$ws.Range("L58").Value = 6656.125
$ws.Range("M58").Value = -1629.2307
$ws.Range("N58").Value = -7062.125
$ws.Range("H132").Value = 2250.4707
$ws.Range("I132").Value = 1323.2222
$ws.Range("J132").Value = 3293.625
$ws.Range("K132").Value = 3969.6666
$ws.Range("L132").Value = 9880.875
$ws.Range("M132").Value = -1439.6666
$ws.Range("N132").Value = -14940.875
$ws.Range("H134").Value = 2409.8333
$ws.Range("I134").Value = 1157.091
$ws.Range("J134").Value = 4378.4287
$ws.Range("K134").Value = 3471.273
$ws.Range("L134").Value = 13135.2861
$ws.Range("M134").Value = -936.2729999999997
$ws.Range("N134").Value = -18205.2861
$ws.Range("H136").Value = 2967.2646
$ws.Range("I136").Value = 1832.2307
$ws.Range("J136").Value = 6656.125
$ws.Range("K136").Value = 5496.6921
$ws.Range("L136").Value = 19968.375
$ws.Range("M136").Value = -2946.6921
$ws.Range("N136").Value = -25068.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 508.5
$ws.Range("I15").Value = 294.66666
$ws.Range("K15").Value = 883.9999799999999
$ws.Range("M15").Value = -743.9999799999999
$ws.Range("H122").Value = 2979.611
$ws.Range("J122").Value = 3386.2258
$ws.Range("L122").Value = 30476.0322
$ws.Range("N122").Value = -35376.0322
$ws.Range("H131").Value = 239314.62
$ws.Range("I131").Value = 3333806.8
$ws.Range("J131").Value = 1276.7693
$ws.Range("K131").Value = 10001420.4
$ws.Range("L131").Value = 3830.3079
$ws.Range("M131").Value = -9996380.399999999
$ws.Range("N131").Value = -13910.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2306.25
$ws.Range("I113").Value = 2421.4285
$ws.Range("K113").Value = 2421.4285
$ws.Range("M113").Value = -251.4285
$ws.Range("H132").Value = 2789.5
$ws.Range("I132").Value = 2959
$ws.Range("K132").Value = 8877
$ws.Range("M132").Value = -6347

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1150.9166
$ws.Range("J22").Value = 1113.3334
$ws.Range("L22").Value = 1113.3334
$ws.Range("N22").Value = -1703.3334
$ws.Range("H27").Value = 1150.9166
$ws.Range("J27").Value = 1113.3334
$ws.Range("L27").Value = 1113.3334
$ws.Range("N27").Value = -1327.3334
$ws.Range("H122").Value = 6560.143
$ws.Range("I122").Value = 9726.933999999999
$ws.Range("J122").Value = 2906.1538
$ws.Range("K122").Value = 29180.802
$ws.Range("L122").Value = 8718.4614
$ws.Range("M122").Value = -26730.802
$ws.Range("N122").Value = -13618.4614
$ws.Range("H132").Value = 66670420
$ws.Range("I132").Value = 90912300
$ws.Range("J132").Value = 5248
$ws.Range("K132").Value = 272736900
$ws.Range("L132").Value = 15744
$ws.Range("M132").Value = -272734370
$ws.Range("N132").Value = -20804
$ws.Range("H136").Value = 55558800
$ws.Range("I136").Value = 83334200
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 250002600
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -250000050
$ws.Range("N136").Value = -29100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2735.2593
$ws.Range("I132").Value = 2321.6365
$ws.Range("J132").Value = 3019.625
$ws.Range("K132").Value = 6964.9095
$ws.Range("L132").Value = 9058.875
$ws.Range("M132").Value = -4434.9095
$ws.Range("N132").Value = -14118.875
$ws.Range("H136").Value = 7577432.5
$ws.Range("I136").Value = 19232034
$ws.Range("J136").Value = 1941.75
$ws.Range("K136").Value = 57696102
$ws.Range("L136").Value = 5825.25
$ws.Range("M136").Value = -57693552
$ws.Range("N136").Value = -10925.25
